# Benchmark for V12 (pellet direction)
# Adds a new row (14) to Sheet1 with the V12 benchmark results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data
$ws.Range("A14").Value = "v12-22700"
$ws.Range("B14").Value = 908
$ws.Range("C14").Value = "Added pellet direction flag"
$ws.Range("D14").Value = 1086
$ws.Range("E14").Value = 0.18
$ws.Range("F14").Value = 0.37
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 2505.66
$ws.Range("J14").Value = 1280
$ws.Range("K14").Value = 4200
$ws.Range("L14").Value = 221.54
$ws.Range("M14").Value = 114
$ws.Range("N14").Value = 244
$ws.Range("O14").Value = 0.65
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 132.73
$ws.Range("S14").Value = 58.9
$ws.Range("T14").Value = 308.7

# Match the number formats/styles used by the rest of the column (same as row 13)
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("E14").NumberFormat = $ws.Range("E13").NumberFormat
$ws.Range("F14").NumberFormat = $ws.Range("F13").NumberFormat
$ws.Range("G14").NumberFormat = $ws.Range("G13").NumberFormat
$ws.Range("H14").NumberFormat = $ws.Range("H13").NumberFormat
$ws.Range("I14").NumberFormat = $ws.Range("I13").NumberFormat
$ws.Range("J14").NumberFormat = $ws.Range("J13").NumberFormat
$ws.Range("K14").NumberFormat = $ws.Range("K13").NumberFormat
$ws.Range("L14").NumberFormat = $ws.Range("L13").NumberFormat
$ws.Range("M14").NumberFormat = $ws.Range("M13").NumberFormat
$ws.Range("N14").NumberFormat = $ws.Range("N13").NumberFormat
$ws.Range("O14").NumberFormat = $ws.Range("O13").NumberFormat
$ws.Range("P14").NumberFormat = $ws.Range("P13").NumberFormat
$ws.Range("Q14").NumberFormat = $ws.Range("Q13").NumberFormat
$ws.Range("R14").NumberFormat = $ws.Range("R13").NumberFormat
$ws.Range("S14").NumberFormat = $ws.Range("S13").NumberFormat
$ws.Range("T14").NumberFormat = $ws.Range("T13").NumberFormat

# Column E width becomes an explicit custom width (same value, just marked customWidth)
$ws.Columns("E").ColumnWidth = $ws.Columns("E").ColumnWidth

# Update the selection to match the new active cell
$ws.Range("F12").Select()
